$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells take a new value that would otherwise be auto-
# parsed as a number by Excel; the source data stores every Price
# cell as text, so force Text format first (cell-by-cell, since a
# single multi-area Range().NumberFormat assignment only takes on
# the first area) to preserve the original text representation.
foreach ($r in @(4,5,6,8,9,10,11,13,14,15,16,19,20,22,25,26,27,28,29,30,32,33,34,35,38,39,41,42,43,44,45,47,48,49,51)) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "36.888.83"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "1.992.25"
$ws.Range("E3").Value = "  -3.41%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.61%  "

$ws.Range("D5").Value = "224.08"
$ws.Range("E5").Value = "  -3.35%  "

$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -3.29%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "54.01"
$ws.Range("E8").Value = "  -5.96%  "

$ws.Range("D9").Value = "0.373"
$ws.Range("E9").Value = "  -3.39%  "

$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -4.29%  "

$ws.Range("D12").Value = "2.293.72"
$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("D13").Value = "13.96"
$ws.Range("E13").Value = "  -5.00%  "

$ws.Range("D14").Value = "19.92"
$ws.Range("E14").Value = "  -5.87%  "

$ws.Range("D15").Value = "0.731"
$ws.Range("E15").Value = "  -4.03%  "

$ws.Range("D16").Value = "5.05"
$ws.Range("E16").Value = "  -5.60%  "

$ws.Range("D17").Value = "1.974.96"
$ws.Range("E17").Value = "  -4.96%  "

$ws.Range("D18").Value = "36.844.60"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("D20").Value = "68.32"
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("D21").Value = "0.0₃0806"
$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("D22").Value = "221.59"
$ws.Range("E22").Value = "  -2.47%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  -7.38%  "

$ws.Range("D26").Value = "165.16"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").Value = "9.02"
$ws.Range("E27").Value = "  -9.70%  "

$ws.Range("D28").Value = "0.123"
$ws.Range("E28").Value = "  -5.86%  "

$ws.Range("D29").Value = "18.52"
$ws.Range("E29").Value = "  -3.80%  "

$ws.Range("D30").Value = "1.33"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("E31").Value = "  -4.78%  "

$ws.Range("D32").Value = "4.45"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("D33").Value = "0.0604"
$ws.Range("E33").Value = "  -3.19%  "

$ws.Range("D34").Value = "4.36"
$ws.Range("E34").Value = "  -5.53%  "

$ws.Range("D35").Value = "2.30"
$ws.Range("E35").Value = "  -8.08%  "

$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  -5.66%  "

$ws.Range("D39").Value = "5.21"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").Value = "1.457.04"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").Value = "0.0213"
$ws.Range("E41").Value = "  -5.69%  "

$ws.Range("D42").Value = "94.05"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0912"
$ws.Range("E43").Value = "  -4.97%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "16.17"
$ws.Range("E44").Value = "  -3.87%  "

$ws.Range("D45").Value = "2.74"
$ws.Range("E45").Value = "  -5.13%  "

$ws.Range("E46").Value = "  -7.39%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "7.10"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.997"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("D50").Value = "2.181.52"
$ws.Range("E50").Value = "  -3.03%  "

$ws.Range("D51").Value = "43.92"
$ws.Range("E51").Value = "  -4.39%  "
